# "suite et correction acp"
# Adds a "Moyenne " (average) summary row beneath the temperature table:
# row 31, column A holds the label, columns B:M hold =AVERAGE() of the
# 28 city rows (2:29) for that month, entered as one shared formula.
# Row 30 is intentionally left blank (matches the source layout).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A31").Value = "Moyenne "

$ws.Range("B31").Formula = "=AVERAGE(B2:B29)"
$ws.Range("C31:M31").Formula = "=AVERAGE(C2:C29)"

# Restore the cursor/selection to where the author left off after typing
# the new row.
$ws.Range("F33").Select()
